# Update simulation results for the 380 kV case (loading_percent sheet).
# Columns B,C,D,F,G,I,J,K,L,M,N,O are refreshed for rows 2-25 (row index col A, header row 1 untouched).
# Columns E and H remain 0 and are left unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "C", "D", "F", "G", "I", "J", "K", "L", "M", "N", "O")

$newValues = @{
    2 = @(13.23879683272692, 7.374573703866711, 8.565971119110056, 36.90773688216225, 3.68644197314729, 27.72940133636863, 10.64593046196853, 10.40890072525362, 11.76800971398558, 15.67100367630393, 21.27667634020186, 28.21702627328564)
    3 = @(13.06064626348497, 7.327300951628263, 8.557519310137183, 36.99010089706999, 3.688087619149918, 27.81445593785877, 10.6636394059248, 10.2798133488081, 11.77691754203949, 15.64926580040314, 21.33392700332621, 28.29504439168336)
    4 = @(12.95224757951528, 7.297678342755244, 8.553421971404388, 37.04718640763265, 3.689152808600792, 27.87066654428886, 10.67513072527262, 10.20132677607724, 11.78367951472343, 15.6378499011456, 21.37074810485476, 28.34745305074139)
    5 = @(12.9083735340066, 7.285460430676398, 8.552028649505186, 37.07208624985042, 3.689600692294563, 27.8945757584945, 10.67996933773084, 10.16957225743996, 11.78676083134633, 15.63368729709407, 21.38617387230296, 28.36994259799803)
    6 = @(12.90110785063717, 7.283422916789102, 8.551814030861788, 37.0763196930441, 3.68967589837171, 27.89860645316348, 10.68078220824033, 10.16431431183621, 11.78729217739507, 15.63302576978566, 21.38876076390943, 28.37374536714494)
    7 = @(12.95165459946572, 7.297514155895465, 8.553402059301863, 37.04751558862992, 3.689158792944193, 27.87098493087813, 10.6751953490815, 10.20089754887175, 11.78371975050514, 15.63779177608749, 21.3709544364773, 28.3477517672131)
    8 = @(13.177196573346, 7.358400135126574, 8.562831358503169, 36.9347836037031, 3.686998052243682, 27.75790098095894, 10.65190849995347, 10.36425183212272, 11.77081332425482, 15.66311000147417, 21.29607069837157, 28.24299170565647)
    9 = @(13.62505672779766, 7.472903618635578, 8.589905989589344, 36.76543375648259, 3.683193420885196, 27.56775689692912, 10.61112731320581, 10.68918778067135, 11.75572607461369, 15.72790496487076, 21.16241091482969, 28.07331614942785)
    10 = @(13.95449890314366, 7.553840023270805, 8.614914926030689, 36.67257918522053, 3.680659222133679, 27.44730290120589, 10.58411649335439, 10.92866541028794, 11.75082692638755, 15.78449237320316, 21.07217125151347, 27.97046957276508)
    11 = @(14.10383584350731, 7.589924474120552, 8.627376091310403, 36.63719580913988, 3.679562475136226, 27.39667846709832, 10.57246377191405, 11.03734286119704, 11.74993051577171, 15.81212746904997, 21.0328305250462, 27.92842111944178)
    12 = @(14.1602595802339, 7.60347951639636, 8.632248227921202, 36.62478294950483, 3.679155186291641, 27.37810769136532, 10.56814202669695, 11.07842345247748, 11.74978163245942, 15.82285901650707, 21.01817783479503, 27.91317979277624)
    13 = @(14.14811415932393, 7.600565119445133, 8.631192147661917, 36.62741242713152, 3.679242546895055, 27.38208057509939, 10.56906875550098, 11.06957982333111, 11.7498052381731, 15.82053601270828, 21.0213226836902, 27.9164319700328)
    14 = @(14.10848068859233, 7.591041868968394, 8.627773866511465, 36.63615483407542, 3.679528806610863, 27.39513862214974, 10.57210639993316, 11.04072425335064, 11.74991445373844, 15.81300504648234, 21.03162013955947, 27.92715354610488)
    15 = @(14.08418602714176, 7.585194249159715, 8.625699964352872, 36.64163822265829, 3.679705192986951, 27.40321513432398, 10.57397887010434, 11.02303881233065, 11.75000613822882, 15.80842668181356, 21.03795947577438, 27.93380958286505)
    16 = @(13.94472436468615, 7.551466723607398, 8.614122164055816, 36.67502955767324, 3.680732022221574, 27.45069523982927, 10.58489076468638, 10.9215547531541, 11.75091223215089, 15.78272398569391, 21.07477657325528, 27.97331288832826)
    17 = @(13.85899716936079, 7.530585271781189, 8.607295520639404, 36.69727034482803, 3.681376283303429, 27.48089085905978, 10.59174713793761, 10.85920489452259, 11.75180861024522, 15.76743726810296, 21.09779981256007, 27.99876039765378)
    18 = @(13.809642423013, 7.518505958105896, 8.603471254278533, 36.71070804680741, 3.681752125303038, 27.4986511541308, 10.59575049946264, 10.8233203060296, 11.75244964317237, 15.75882329401849, 21.11120316386195, 28.01384298115069)
    19 = @(13.79292530643378, 7.514404401254413, 8.602194059583367, 36.71536866155358, 3.681880286958688, 27.50473191415147, 10.59711624361282, 10.81116761584196, 11.7526882635545, 15.75593757887454, 21.11576899182035, 28.01902624855119)
    20 = @(13.8681281864087, 7.532815292008219, 8.60801166551709, 36.69483597830276, 3.681307154418222, 27.47763585499127, 10.59101108338803, 10.86584474235263, 11.75170021175419, 15.76904612602066, 21.09533229504642, 27.99600532048251)
    21 = @(14.12012584304439, 7.593842076568929, 8.628773757487847, 36.63356021651998, 3.67944450768417, 27.39128688537353, 10.57121170641887, 11.04920209881352, 11.7498772108031, 15.81520987969312, 21.02858889213319, 27.92398586090604)
    22 = @(14.28405856535726, 7.633087460561991, 8.643235622248701, 36.59926031924137, 3.678273921496053, 27.33834769909974, 10.55880126070172, 11.16859426630466, 11.74979604560694, 15.84693258826565, 20.98639454858802, 27.88088945983171)
    23 = @(14.19665077627052, 7.612201191454475, 8.635436261184655, 36.61704097428174, 3.678894419114861, 27.36628260544147, 10.56537661354609, 11.10492439249178, 11.74973812126236, 15.82986144962641, 21.00878430572695, 27.90352724063835)
    24 = @(13.86400026474054, 7.531807330982063, 8.607687583307547, 36.69593452705639, 3.68133839061557, 27.47910619500324, 10.59134366192804, 10.86284298535835, 11.75174882711991, 15.76831821719685, 21.09644733948203, 27.99724948147938)
    25 = @(13.5036139971973, 7.442472654285064, 8.581674060748098, 36.80570585527993, 3.684176639046539, 27.61581435705477, 10.62163956007221, 10.601005059773, 11.75871779793538, 15.70877961185654, 21.19716587059789, 28.11538798870364)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $row).Value = $vals[$i]
    }
}
